$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 185-186; all existing data from row 185
# downward shifts down by two rows (old row 185 -> 187, old row 186 -> 188,
# ..., old row 293 -> 295, old row 294 -> 296).
$ws.Rows("185:186").Insert()

# Seed the two new rows with a copy of the (now shifted) rows that used to
# sit directly below them, then overwrite just the cells that actually hold
# new data, per the new weekly price report entries.
$ws.Range("A187:T187").Copy()
$ws.Range("A185").PasteSpecial()

$ws.Range("A188:T188").Copy()
$ws.Range("A186").PasteSpecial()

# New row 185: Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg
$ws.Range("D185").Value = 44719
$ws.Range("M185").Value = 160
$ws.Range("N185").Value = 14000
$ws.Range("O185").Value = 15000
$ws.Range("P185").Value = 14500
$ws.Range("S185").Value = 1036

# New row 186: Fecha, Calidad, Precio minimo/maximo/promedio, Precio $/Kg
$ws.Range("D186").Value = 44719
$ws.Range("L186").Value = "Segunda"
$ws.Range("N186").Value = 12000
$ws.Range("O186").Value = 12000
$ws.Range("P186").Value = 12000
$ws.Range("S186").Value = 857
